$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value (values are stored as text/inline strings in the sheet)
$updates = @{
    "D2"  = "243.61"
    "D3"  = "23.00"
    "D4"  = "5.419"
    "D5"  = "0.05948"
    "D6"  = "3.450"
    "D7"  = "6.555"
    "D8"  = "0.8124"
    "D9"  = "0.9180"
    "D11" = "0.07491"
    "D12" = "0.03273"
    "D13" = "0.03057"
    "D15" = "3.862"
    "D16" = "0.001572"
    "D17" = "0.04670"
    "D18" = "0.0005937"
    "D19" = "0.006048"
    "D20" = "0.004991"
    "D21" = "0.0009806"
    "D22" = "0.0001100"
    "E22" = "21NitroExNTX"
    "D26" = "0.1295"
    "D40" = "0.03948"
    "D41" = "0.006194"
    "E41" = "40KickTokenKICKBestin24h"
    "D43" = "0.002999"
    "D44" = "0.008089"
    "D45" = "0.00005201"
    "D48" = "0.8996"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $updates[$cell]
}
